# Update "想去人数" (number of people interested) counts that changed between
# crawl runs, as reflected in the regenerated gh-pages data output.
#
# The same set of events is listed both on the "展览" sheet and on the
# aggregated "全部类型" sheet, so each changed value needs to be updated in
# both places (rows differ between the two sheets because "全部类型" also
# interleaves rows from the other category sheets).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 9807
$wsExhibition.Range("F9").Value = 716
$wsExhibition.Range("F13").Value = 3037
$wsExhibition.Range("F16").Value = 1993
$wsExhibition.Range("F20").Value = 1569
$wsExhibition.Range("F23").Value = 201
$wsExhibition.Range("F26").Value = 354
$wsExhibition.Range("F28").Value = 332
$wsExhibition.Range("F33").Value = 227
$wsExhibition.Range("F35").Value = 75
$wsExhibition.Range("F36").Value = 377

# Sheet "全部类型" (All types - aggregated view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 9807
$wsAll.Range("F11").Value = 716
$wsAll.Range("F15").Value = 3037
$wsAll.Range("F18").Value = 1993
$wsAll.Range("F22").Value = 1569
$wsAll.Range("F25").Value = 201
$wsAll.Range("F28").Value = 354
$wsAll.Range("F30").Value = 332
$wsAll.Range("F39").Value = 227
$wsAll.Range("F41").Value = 75
$wsAll.Range("F43").Value = 377
